$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# NOTE on ordering: the shared-strings table records new unique strings in the
# order cells are first written. The target file's sharedStrings table has
# "John Doe" / "john.doe@xyz.com" inserted before "Jane Smith" /
# "jane.smith@xyz.com", even though John Doe ends up in row 32 (below Jane
# Smith in row 31). So we populate row 32's text cells before row 31's.

# Row 32 text first (John Doe) so its strings land earlier in sharedStrings.xml
$ws.Cells.Item(32, 3).Value = "John Doe"
$ws.Cells.Item(32, 4).Value = "john.doe@xyz.com"

# Row 31 text next (Jane Smith)
$ws.Cells.Item(31, 3).Value = "Jane Smith"
$ws.Cells.Item(31, 4).Value = "jane.smith@xyz.com"

# Row 31: Jane Smith - remaining columns
$ws.Cells.Item(31, 1).Value = 110030
$ws.Cells.Item(31, 2).Value = 9317596768
$ws.Cells.Item(31, 5).Value = 818876432
$ws.Cells.Item(31, 6).Value = "ACT"
$ws.Cells.Item(31, 7).Value = "eng"
$ws.Cells.Item(31, 8).Value = "PWD"
$ws.Cells.Item(31, 9).Value = $true
$ws.Cells.Item(31, 9).HorizontalAlignment = -4131
$ws.Cells.Item(31, 10).Value = "superadmin"
$ws.Cells.Item(31, 11).Value = "now()"
$ws.Cells.Item(31, 12).Value = "now()"

# Row 32: John Doe - remaining columns
$ws.Cells.Item(32, 1).Value = 110031
$ws.Cells.Item(32, 2).Value = 9317596767
$ws.Cells.Item(32, 5).Value = 818876431
$ws.Cells.Item(32, 6).Value = "ACT"
$ws.Cells.Item(32, 7).Value = "eng"
$ws.Cells.Item(32, 8).Value = "PWD"
$ws.Cells.Item(32, 9).Value = $true
$ws.Cells.Item(32, 9).HorizontalAlignment = -4131
$ws.Cells.Item(32, 10).Value = "superadmin"
$ws.Cells.Item(32, 11).Value = "now()"
$ws.Cells.Item(32, 12).Value = "now()"

$ws.Range("F30").Select()
